$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCD0321")

# Split the combined query into two separate query cells, and duplicate the
# explanation text into the second explanation column.
$ws.Range("T2").Value = "SELECT * FROM DigisalesNew..Tbl_Master_Role WHERE Id=1005 "
$ws.Range("U2").Value = "SELECT * FROM DigisalesNew..Tbl_Pegawai WHERE Role_Id=1005"
$ws.Range("W2").Value = "Melakukan melakuakn syncrnz dan data sesuai"
$ws.Range("X2").Value = "Melakukan melakuakn syncrnz dan data sesuai"

# Widen column U to fit the new content
$ws.Columns.Item(21).ColumnWidth = 16.85546875

# Update selection / view
$ws.Range("X2").Select()
